{"js": "// Fix #148: Endusercontent can be missing in the generated document.\n//\n// The paragraph immediately preceding the `m:endusercontent` fldSimple\n// field used to carry a leftover `_GoBack` bookmark (bookmarkStart /\n// bookmarkEnd). That bookmark is removed and replaced with three\n// spellStart/spellEnd `w:proofErr` marker pairs (matching the proofing\n// marks already used around the \"Some protected text.\" paragraph just\n// above it).\n\n// 1. Remove the stray \"_GoBack\" bookmark.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2. Locate the paragraph that holds the `m:endusercontent` field \u2014 the\n//    paragraph that used to contain the bookmark (4th paragraph, 0-based\n//    index 3: \"A simple demonstration\u2026\", \"\" (usercontent field), \"Some\n//    protected text.\", \"\" (endusercontent field), \"End of demonstration.\",\n//    \"\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[3];\n\n// 3. Insert three `w:proofErr` spellStart/spellEnd pairs at the very\n//    start of that paragraph (i.e. right before the `m:endusercontent`\n//    fldSimple), mirroring the markup produced by Word's proofing pass.\nconst proofErrOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"spellEnd\"/>' +\n  '</w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntargetParagraph.insertOoxml(proofErrOoxml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Fixed #148 Endusercontent can be missing in the generated document.\n#\n# The paragraph immediately preceding the `m:endusercontent` fldSimple\n# field used to carry a leftover \"_GoBack\" bookmark (bookmarkStart /\n# bookmarkEnd). That bookmark is removed and replaced with three\n# spellStart/spellEnd `w:proofErr` marker pairs (matching the proofing\n# marks already used around the \"Some protected text.\" paragraph just\n# above it).\n\n$d = $word.ActiveDocument\n\n# 1. Remove the stray \"_GoBack\" bookmark (hidden bookmark, but still\n#    reachable by name through the Bookmarks collection).\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bm.Delete()\n\n# 2. Build the OOXML fragment for the three proofErr spellStart/spellEnd\n#    pairs that replace the bookmark.\n$proofErrXml = @'\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"spellEnd\"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n\n# 3. Insert it right at the boundary between the \"Some protected text.\"\n#    paragraph and the (now bookmark-free) paragraph that holds the\n#    `m:endusercontent` field \u2014 i.e. a zero-length range positioned at\n#    the end of paragraph 3 / start of paragraph 4 \u2014 so the three\n#    proofErr pairs land at the very beginning of the target paragraph,\n#    right before the `m:endusercontent` fldSimple, leaving that field\n#    untouched.\n$precedingParagraph = $d.Paragraphs.Item(3)\n$insertionPoint = $d.Range($precedingParagraph.Range.End, $precedingParagraph.Range.End)\n$insertionPoint.InsertXML($proofErrXml)\n"}
